$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Añadida hora estimada hacer test:
# new row 15 -> task "Hacer test", 5 estimated hours, owner "Ángel"
$ws.Range("A15").Value = "Hacer test"
$ws.Range("B15").Value = 5
$ws.Range("D15").Value = "Ángel"

# Leave the selection where the author ended up after entering the row
$ws.Range("B17").Select()
